$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1424151111754526
$ws.Range("D2").Value = 0.01443668439061518
$ws.Range("E2").Value = 0.4178968709247641
$ws.Range("F2").Value = 0.8588185927887366
$ws.Range("G2").Value = 0.7181621447065254
$ws.Range("H2").Value = 0.7488212154705138
$ws.Range("K2").Value = 0.9119527343079596

$ws.Range("B3").Value = 0.1329425658592243
$ws.Range("D3").Value = 0.01451165562322032
$ws.Range("E3").Value = 0.3643918641718642
$ws.Range("F3").Value = 0.838997368744117
$ws.Range("G3").Value = 0.697448958700349
$ws.Range("H3").Value = 0.7452179121359848
$ws.Range("K3").Value = 0.7949186235302363

$ws.Range("B4").Value = 0.1271988809211138
$ws.Range("D4").Value = 0.01456171990234267
$ws.Range("E4").Value = 0.3316558605631457
$ws.Range("F4").Value = 0.8275570947214419
$ws.Range("G4").Value = 0.6854153286938782
$ws.Range("H4").Value = 0.7435364794493609
$ws.Range("K4").Value = 0.7229519513519165

$ws.Range("B5").Value = 0.1248766121654654
$ws.Range("D5").Value = 0.01458313608027773
$ws.Range("E5").Value = 0.3183424719556029
$ws.Range("F5").Value = 0.8230772890839688
$ws.Range("G5").Value = 0.6806819383414506
$ws.Range("H5").Value = 0.7429841627452163
$ws.Range("K5").Value = 0.6935976462954443

$ws.Range("B6").Value = 0.1244921112266155
$ws.Range("D6").Value = 0.01458675353427097
$ws.Range("E6").Value = 0.3161333415720549
$ws.Range("F6").Value = 0.8223443850476428
$ws.Range("G6").Value = 0.6799062060889014
$ws.Range("H6").Value = 0.7429004578237226
$ws.Range("K6").Value = 0.6887217285355405

$ws.Range("B7").Value = 0.1271674876436606
$ws.Range("D7").Value = 0.01456200461994328
$ws.Range("E7").Value = 0.3314762062300503
$ws.Range("F7").Value = 0.8274959426062054
$ws.Range("G7").Value = 0.6853508047655197
$ws.Range("H7").Value = 0.7435284935495616
$ws.Range("K7").Value = 0.7225561798745161

$ws.Range("B8").Value = 0.1391339872389921
$ws.Range("D8").Value = 0.0144616979586516
$ws.Range("E8").Value = 0.3994227176412579
$ws.Range("F8").Value = 0.8518318779889427
$ws.Range("G8").Value = 0.7108771843748087
$ws.Range("H8").Value = 0.7474681922807633
$ws.Range("K8").Value = 0.8716213881060355

$ws.Range("B9").Value = 0.1631724843750817
$ws.Range("D9").Value = 0.01429699053659128
$ws.Range("E9").Value = 0.5337051909546773
$ws.Range("F9").Value = 0.9054127890319421
$ws.Range("G9").Value = 0.7664433230853547
$ws.Range("H9").Value = 0.7594379314803916
$ws.Range("K9").Value = 1.163117941388009

$ws.Range("B10").Value = 0.1811804329263538
$ws.Range("D10").Value = 0.01419554867388939
$ws.Range("E10").Value = 0.6331678527653679
$ws.Range("F10").Value = 0.9484462813745012
$ws.Range("G10").Value = 0.8107405677280894
$ws.Range("H10").Value = 0.7708640545827166
$ws.Range("K10").Value = 1.37684716269365

$ws.Range("B11").Value = 0.1894477564342196
$ws.Range("D11").Value = 0.01415367986103533
$ws.Range("E11").Value = 0.6786289171441666
$ws.Range("F11").Value = 0.9688409186484961
$ws.Range("G11").Value = 0.8316717099166908
$ws.Range("H11").Value = 0.7766434580529449
$ws.Range("K11").Value = 1.474000398213832

$ws.Range("B12").Value = 0.1925891547395793
$ws.Range("D12").Value = 0.01413844401657016
$ws.Range("E12").Value = 0.6958777122517716
$ws.Range("F12").Value = 0.9766831656735491
$ws.Range("G12").Value = 0.839711919813368
$ws.Range("H12").Value = 0.7789163732803388
$ws.Range("K12").Value = 1.510780179962978

$ws.Range("B13").Value = 0.1919121222641422
$ws.Range("D13").Value = 0.0141416977291513
$ws.Range("E13").Value = 0.6921613334265402
$ws.Range("F13").Value = 0.9749888694867224
$ws.Range("G13").Value = 0.8379752149447484
$ws.Range("H13").Value = 0.7784230954931388
$ws.Range("K13").Value = 1.502859434076015

$ws.Range("B14").Value = 0.1897059864550101
$ws.Range("D14").Value = 0.01415241396862044
$ws.Range("E14").Value = 0.6800472950032201
$ws.Range("F14").Value = 0.9694837063866544
$ws.Range("G14").Value = 0.8323308865575427
$ws.Range("H14").Value = 0.7768287567407697
$ws.Range("K14").Value = 1.477026491220329

$ws.Range("B15").Value = 0.1883560605112393
$ws.Range("D15").Value = 0.01415905871039413
$ws.Range("E15").Value = 0.6726315650329298
$ws.Range("F15").Value = 0.966127209693056
$ws.Range("G15").Value = 0.8288884794347382
$ws.Range("H15").Value = 0.7758631893521795
$ws.Range("K15").Value = 1.461201775737038

$ws.Range("B16").Value = 0.1806416526886494
$ws.Range("D16").Value = 0.0141983712867102
$ws.Range("E16").Value = 0.6302014265967273
$ws.Range("F16").Value = 0.947130045510022
$ws.Range("G16").Value = 0.8093885249269874
$ws.Range("H16").Value = 0.770498123972601
$ws.Range("K16").Value = 1.370496537373413

$ws.Range("B17").Value = 0.1759283562378045
$ws.Range("D17").Value = 0.01422358635643661
$ws.Range("E17").Value = 0.6042289279725424
$ws.Range("F17").Value = 0.9356866283153806
$ws.Range("G17").Value = 0.7976270583305052
$ws.Range("H17").Value = 0.7673563215049057
$ws.Range("K17").Value = 1.314833386063185

$ws.Range("B18").Value = 0.1732244974942603
$ws.Range("D18").Value = 0.01423849174496183
$ws.Range("E18").Value = 0.5893103029655151
$ws.Range("F18").Value = 0.929181624126258
$ws.Range("G18").Value = 0.7909354886209599
$ws.Range("H18").Value = 0.7656039379382094
$ws.Range("K18").Value = 1.282810485814878

$ws.Range("B19").Value = 0.1723102407663504
$ws.Range("D19").Value = 0.01424360745325615
$ws.Range("E19").Value = 0.5842624764112117
$ws.Range("F19").Value = 0.9269923103352085
$ws.Range("G19").Value = 0.7886823734580162
$ws.Range("H19").Value = 0.7650199827868391
$ws.Range("K19").Value = 1.271966877819182

$ws.Range("B20").Value = 0.1764293601259084
$ws.Range("D20").Value = 0.01422086050396487
$ws.Range("E20").Value = 0.6069916477884618
$ws.Range("F20").Value = 0.9368968248045775
$ws.Range("G20").Value = 0.7988714852960186
$ws.Range("H20").Value = 0.7676851053350617
$ws.Range("K20").Value = 1.320759537080676

$ws.Range("B21").Value = 0.1903536912892605
$ws.Range("D21").Value = 0.01414924951573582
$ws.Range("E21").Value = 0.683604545553905
$ws.Range("F21").Value = 0.9710974577167519
$ws.Range("G21").Value = 0.8339856526763754
$ws.Range("H21").Value = 0.77729475662494
$ws.Range("K21").Value = 1.484614521059825

$ws.Range("B22").Value = 0.1995166141811922
$ws.Range("D22").Value = 0.0141060577748533
$ws.Range("E22").Value = 0.7338731084087584
$ws.Range("F22").Value = 0.9941452992306949
$ws.Range("G22").Value = 0.8576003821598874
$ws.Range("H22").Value = 0.7840674058667219
$ws.Range("K22").Value = 1.591644965317641

$ws.Range("B23").Value = 0.1946204984979829
$ws.Range("D23").Value = 0.01412877820343894
$ws.Range("E23").Value = 0.7070248262493806
$ws.Range("F23").Value = 0.9817800577671392
$ws.Range("G23").Value = 0.8449352355766848
$ws.Range("H23").Value = 0.7804074316231322
$ws.Range("K23").Value = 1.534525958759843

$ws.Range("B24").Value = 0.1762028379515357
$ws.Range("D24").Value = 0.01422209158793208
$ws.Range("E24").Value = 0.6057425809381272
$ws.Range("F24").Value = 0.9363494646806743
$ws.Range("G24").Value = 0.7983086611712054
$ws.Range("H24").Value = 0.7675362943824382
$ws.Range("K24").Value = 1.318080390953071

$ws.Range("B25").Value = 0.1566083816444888
$ws.Range("D25").Value = 0.01433811940731466
$ws.Range("E25").Value = 0.4972493349182088
$ws.Range("F25").Value = 0.8902806614828052
$ws.Range("G25").Value = 0.7508096592722779
$ws.Range("H25").Value = 0.7557408544608109
$ws.Range("K25").Value = 1.084341454631442

